# Weekly update: insert a new price record as row 255, shifting the
# existing rows 255-279 down to 256-280.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 255 (this shifts rows 255..279
# down to 256..280 and copies formatting from the row above, matching the
# date style used throughout column D).
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new weekly data point.
$ws.Cells.Item(255, 1).Value = 8
$ws.Cells.Item(255, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(255, 3).Value = "Coquimbo"
$ws.Cells.Item(255, 4).Value = 44946
$ws.Cells.Item(255, 5).Value = 4
$ws.Cells.Item(255, 6).Value = 100112037
$ws.Cells.Item(255, 7).Value = "Cebollín"
$ws.Cells.Item(255, 8).Value = "Sin especificar"
$ws.Cells.Item(255, 9).Value = "Primera"
$ws.Cells.Item(255, 10).Value = 1000
$ws.Cells.Item(255, 11).Value = 1200
$ws.Cells.Item(255, 12).Value = 1400
$ws.Cells.Item(255, 13).Value = 1300
$ws.Cells.Item(255, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(255, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(255, 16).Value = 217
$ws.Cells.Item(255, 17).Value = 6
$ws.Cells.Item(255, 18).Value = "Hortaliza"
